# Automated monthly rollover update
# - "VENTAS POR GRUPO": zero out a handful of now-stale figures (and their
#   matching "N de 33" progress labels in the totals row).
# - "VENTA MENSUAL": shift the rolling 4-month window one month forward
#   (junio-septiembre -> julio-octubre), which shifts each advisor's C/D/E/F
#   figures one column to the left and drops in 0 for the new (not yet
#   reported) rightmost month; also re-narrows the month columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("M2").Value = 0

$wsGrupo.Range("L3").Value = 0
$wsGrupo.Range("M3").Value = 0

$wsGrupo.Range("M8").Value = 0

$wsGrupo.Range("M10").Value = 0

$wsGrupo.Range("D12").Value = 0
$wsGrupo.Range("L12").Value = 0

$wsGrupo.Range("D15").Value = 0
$wsGrupo.Range("M15").Value = 0

$wsGrupo.Range("L19").Value = 0
$wsGrupo.Range("M19").Value = 0

$wsGrupo.Range("E20").Value = 0
$wsGrupo.Range("H20").Value = 0
$wsGrupo.Range("M20").Value = 0
$wsGrupo.Range("O20").Value = 0
$wsGrupo.Range("P20").Value = 0

$wsGrupo.Range("E27").Value = 0
$wsGrupo.Range("H27").Value = 0

$wsGrupo.Range("D28").Value = 0

$wsGrupo.Range("M34").Value = 0

# Totals row: progress counters ("N de 33") for the columns touched above
$wsGrupo.Range("D35").Value = "0 de 33"
$wsGrupo.Range("E35").Value = "0 de 33"
$wsGrupo.Range("H35").Value = "0 de 33"
$wsGrupo.Range("L35").Value = "0 de 33"
$wsGrupo.Range("M35").Value = "0 de 33"
$wsGrupo.Range("O35").Value = "0 de 33"
$wsGrupo.Range("P35").Value = "0 de 33"

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Re-narrow/re-widen the month columns (C..F) to their new widths.
# ColumnWidth uses Excel's character-width units, which differ from the
# raw <col width> stored in the file by a constant font padding offset;
# 0.8333333333333339 is that offset for this workbook's default font.
$wsMensual.Columns.Item(3).ColumnWidth = 12.166666666666666
$wsMensual.Columns.Item(4).ColumnWidth = 13.166666666666666
$wsMensual.Columns.Item(5).ColumnWidth = 15.166666666666666
$wsMensual.Columns.Item(6).ColumnWidth = 12.166666666666666

# Headers roll forward by one month
$wsMensual.Range("C1").Value = "julio"
$wsMensual.Range("D1").Value = "agosto"
$wsMensual.Range("E1").Value = "septiembre"
$wsMensual.Range("F1").Value = "octubre"

# Row 2
$wsMensual.Range("C2").Value = 4360.63
$wsMensual.Range("D2").Value = 275.07
$wsMensual.Range("E2").Value = 6231.33
$wsMensual.Range("F2").Value = 0

# Row 3
$wsMensual.Range("E3").Value = 687.03
$wsMensual.Range("F3").Value = 0

# Row 4
$wsMensual.Range("C4").Value = 1190.78
$wsMensual.Range("D4").Value = 0

# Row 6
$wsMensual.Range("D6").Value = 0

# Row 8
$wsMensual.Range("C8").Value = 0
$wsMensual.Range("D8").Value = 2588.16
$wsMensual.Range("E8").Value = 3558.2
$wsMensual.Range("F8").Value = 0

# Row 10
$wsMensual.Range("D10").Value = 3415.45
$wsMensual.Range("E10").Value = 226.8
$wsMensual.Range("F10").Value = 0

# Row 11
$wsMensual.Range("C11").Value = -86.23
$wsMensual.Range("D11").Value = 0

# Row 12
$wsMensual.Range("C12").Value = 112.01
$wsMensual.Range("D12").Value = 1939.26
$wsMensual.Range("E12").Value = 998.71
$wsMensual.Range("F12").Value = 0

# Row 13
$wsMensual.Range("D13").Value = 1831.68
$wsMensual.Range("E13").Value = 0

# Row 15
$wsMensual.Range("C15").Value = 374.03
$wsMensual.Range("D15").Value = 3914.73
$wsMensual.Range("E15").Value = 10853.08
$wsMensual.Range("F15").Value = 0

# Row 19
$wsMensual.Range("E19").Value = 4077.46
$wsMensual.Range("F19").Value = 0

# Row 20
$wsMensual.Range("C20").Value = 5850.44
$wsMensual.Range("D20").Value = 4971.22
$wsMensual.Range("E20").Value = 4354.56
$wsMensual.Range("F20").Value = 0

# Row 23
$wsMensual.Range("C23").Value = 612.28
$wsMensual.Range("D23").Value = 732.5700000000001
$wsMensual.Range("E23").Value = 0

# Row 24
$wsMensual.Range("C24").Value = 128.3
$wsMensual.Range("D24").Value = 0

# Row 25
$wsMensual.Range("D25").Value = 1599.58
$wsMensual.Range("E25").Value = 0

# Row 27
$wsMensual.Range("E27").Value = 948.92
$wsMensual.Range("F27").Value = 0

# Row 28
$wsMensual.Range("C28").Value = 0
$wsMensual.Range("D28").Value = 11570.11
$wsMensual.Range("E28").Value = 1831.68
$wsMensual.Range("F28").Value = 0

# Row 30
$wsMensual.Range("C30").Value = 10.76
$wsMensual.Range("D30").Value = 1093.88
$wsMensual.Range("E30").Value = 0

# Row 32
$wsMensual.Range("D32").Value = 132.58
$wsMensual.Range("E32").Value = 0

# Row 34
$wsMensual.Range("C34").Value = 0
$wsMensual.Range("E34").Value = 5238.25
$wsMensual.Range("F34").Value = 0

# Row 35 (totals)
$wsMensual.Range("C35").Value = 12563.28
$wsMensual.Range("D35").Value = 34064.29
$wsMensual.Range("E35").Value = 39006.02
$wsMensual.Range("F35").Value = 0
